$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9632929301833428
$ws.Range("C2").Value = 0.3291852076740582
$ws.Range("D2").Value = 0.01843078997180925
$ws.Range("F2").Value = 0.5048314631687987
$ws.Range("G2").Value = 0.3485194568486634
$ws.Range("H2").Value = 0.5010858148871762
$ws.Range("L2").Value = 0.3002857918979061
$ws.Range("N2").Value = 1.026827486438648
$ws.Range("O2").Value = 1.635046513147643

$ws.Range("B3").Value = 0.864329848898592
$ws.Range("C3").Value = 0.3251402203457303
$ws.Range("D3").Value = 0.01651906018789617
$ws.Range("F3").Value = 0.4996735839254995
$ws.Range("G3").Value = 0.3446663276027806
$ws.Range("H3").Value = 0.5030186136621495
$ws.Range("L3").Value = 0.2887333619204639
$ws.Range("N3").Value = 1.032039124809209
$ws.Range("O3").Value = 1.630651147356332

$ws.Range("B4").Value = 0.8035916768343441
$ws.Range("C4").Value = 0.322689165744336
$ws.Range("D4").Value = 0.01533791880488167
$ws.Range("F4").Value = 0.4968781587208397
$ws.Range("G4").Value = 0.3426017108718113
$ws.Range("H4").Value = 0.5044746769922952
$ws.Range("L4").Value = 0.2818117439585421
$ws.Range("N4").Value = 1.035661137352122
$ws.Range("O4").Value = 1.629144940626716

$ws.Range("B5").Value = 0.778848423792283
$ws.Range("C5").Value = 0.3216987072246553
$ws.Range("D5").Value = 0.01485477947802849
$ws.Range("F5").Value = 0.4958323704339023
$ws.Range("G5").Value = 0.3418359612611468
$ws.Range("H5").Value = 0.505135772291311
$ws.Range("L5").Value = 0.2790343412250849
$ws.Range("N5").Value = 1.037243417704445
$ws.Range("O5").Value = 1.628830746487651

$ws.Range("B6").Value = 0.7747403590585975
$ws.Range("C6").Value = 0.3215347527210497
$ws.Range("D6").Value = 0.01477444581318821
$ws.Range("F6").Value = 0.4956643558032923
$ws.Range("G6").Value = 0.3417133708038662
$ws.Range("H6").Value = 0.5052496383179559
$ws.Range("L6").Value = 0.2785757672621827
$ws.Range("N6").Value = 1.037512577866629
$ws.Range("O6").Value = 1.628796661437846

$ws.Range("B7").Value = 0.8032579459503779
$ws.Range("C7").Value = 0.3226757739701469
$ws.Range("D7").Value = 0.01533141032862062
$ws.Range("F7").Value = 0.496863676875634
$ws.Range("G7").Value = 0.3425910778032417
$ws.Range("H7").Value = 0.5044833184669599
$ws.Range("L7").Value = 0.2817741118901722
$ws.Range("N7").Value = 1.035682046004723
$ws.Range("O7").Value = 1.629139490614364

$ws.Range("B8").Value = 0.9291660941681243
$ws.Range("C8").Value = 0.3277838351075957
$ws.Range("D8").Value = 0.01777316404363205
$ws.Range("F8").Value = 0.5029758785731318
$ws.Range("G8").Value = 0.3471282898257471
$ws.Range("H8").Value = 0.5016963557535519
$ws.Range("L8").Value = 0.2962669243228362
$ws.Range("N8").Value = 1.028536954120433
$ws.Range("O8").Value = 1.633283297572348

$ws.Range("B9").Value = 1.176215387474201
$ws.Range("C9").Value = 0.3380522585948285
$ws.Range("D9").Value = 0.0225022009267164
$ws.Range("F9").Value = 0.5179137668709757
$ws.Range("G9").Value = 0.3584234119696958
$ws.Range("H9").Value = 0.4983677679927609
$ws.Range("L9").Value = 0.3260483224341328
$ws.Range("N9").Value = 1.017868015395422
$ws.Range("O9").Value = 1.650887646636221

$ws.Range("B10").Value = 1.357748646069581
$ws.Range("C10").Value = 0.3457412348067948
$ws.Range("D10").Value = 0.02593939541286261
$ws.Range("F10").Value = 0.530695703252583
$ws.Range("G10").Value = 0.3681956545490408
$ws.Range("H10").Value = 0.497225059631063
$ws.Range("L10").Value = 0.3487602012798732
$ws.Range("N10").Value = 1.012059826200556
$ws.Range("O10").Value = 1.669626298072956

$ws.Range("B11").Value = 1.440326543711762
$ws.Range("C11").Value = 0.3492688270688546
$ws.Range("D11").Value = 0.02749476047719668
$ws.Range("F11").Value = 0.5369046512184781
$ws.Range("G11").Value = 0.3729639616897771
$ws.Range("H11").Value = 0.496988203940262
$ws.Range("L11").Value = 0.3592735063907355
$ws.Range("N11").Value = 1.009856950946016
$ws.Range("O11").Value = 1.679417373457085

$ws.Range("B12").Value = 1.471594878423957
$ws.Range("C12").Value = 0.3506087621790925
$ws.Range("D12").Value = 0.02808252714331161
$ws.Range("F12").Value = 0.5393126226731795
$ws.Range("G12").Value = 0.3748162031992592
$ws.Range("H12").Value = 0.4969392016114824
$ws.Range("L12").Value = 0.3632807111520009
$ws.Range("N12").Value = 1.009085832146269
$ws.Range("O12").Value = 1.68330754385326

$ws.Range("B13").Value = 1.46486080556906
$ws.Range("C13").Value = 0.3503200031021692
$ws.Range("D13").Value = 0.02795599569116547
$ws.Range("F13").Value = 0.5387914967034959
$ws.Range("G13").Value = 0.374415215355242
$ws.Range("H13").Value = 0.4969479454635319
$ws.Range("L13").Value = 0.3624165306017773
$ws.Range("N13").Value = 1.009249103344956
$ws.Range("O13").Value = 1.68246160378547

$ws.Range("B14").Value = 1.442899061073035
$ws.Range("C14").Value = 0.3493789829349367
$ws.Range("D14").Value = 0.02754314095815857
$ws.Range("F14").Value = 0.5371016181992161
$ws.Range("G14").Value = 0.3731154121319662
$ws.Range("H14").Value = 0.4969833570097535
$ws.Range("L14").Value = 0.3596026599380764
$ws.Range("N14").Value = 1.009792247281794
$ws.Range("O14").Value = 1.67973376054303

$ws.Range("B15").Value = 1.429446529637346
$ws.Range("C15").Value = 0.348803111201434
$ws.Range("D15").Value = 0.02729009614361644
$ws.Range("F15").Value = 0.5360739148730858
$ws.Range("G15").Value = 0.3723253168430176
$ws.Range("H15").Value = 0.4970103465916083
$ws.Range("L15").Value = 0.3578824727447483
$ws.Range("N15").Value = 1.010133148051736
$ws.Range("O15").Value = 1.678086655127203

$ws.Range("B16").Value = 1.352351785204633
$ws.Range("C16").Value = 0.3455112835216596
$ws.Range("D16").Value = 0.02583758021449256
$ws.Range("F16").Value = 0.5302978755471273
$ws.Range("G16").Value = 0.367890545005892
$ws.Range("H16").Value = 0.497246232694522
$ws.Range("L16").Value = 0.348076780061163
$ws.Range("N16").Value = 1.012212619728459
$ws.Range("O16").Value = 1.669011948298419

$ws.Range("B17").Value = 1.305054761948895
$ws.Range("C17").Value = 0.3434993690672883
$ws.Range("D17").Value = 0.02494437716831044
$ws.Range("F17").Value = 0.5268555205513152
$ws.Range("G17").Value = 0.3652527636273248
$ws.Range("H17").Value = 0.497463415663006
$ws.Range("L17").Value = 0.3421077615484336
$ws.Range("N17").Value = 1.013600744707176
$ws.Range("O17").Value = 1.66376960472283

$ws.Range("B18").Value = 1.277850635866628
$ws.Range("C18").Value = 0.3423449882146201
$ws.Range("D18").Value = 0.02442985751956428
$ws.Range("F18").Value = 0.5249126861719233
$ws.Range("G18").Value = 0.3637659616997979
$ws.Range("H18").Value = 0.497614967905804
$ws.Range("L18").Value = 0.3386916328769871
$ws.Range("N18").Value = 1.014440512582496
$ws.Range("O18").Value = 1.660873561392862

$ws.Range("B19").Value = 1.268639817816222
$ws.Range("C19").Value = 0.3419546243277125
$ws.Range("D19").Value = 0.02425551825767513
$ws.Range("F19").Value = 0.5242612486997658
$ws.Range("G19").Value = 0.3632677690469848
$ws.Range("H19").Value = 0.4976708553846123
$ws.Range("L19").Value = 0.3375379280478512
$ws.Range("N19").Value = 1.014731949714701
$ws.Range("O19").Value = 1.659913475498882

$ws.Range("B20").Value = 1.310089636984571
$ws.Range("C20").Value = 0.3437132503210165
$ws.Range("D20").Value = 0.02503954038348155
$ws.Range("F20").Value = 0.5272181231449338
$ws.Range("G20").Value = 0.3655304145650717
$ws.Range("H20").Value = 0.4974375396141681
$ws.Range("L20").Value = 0.3427414051663078
$ws.Range("N20").Value = 1.013448697134308
$ws.Range("O20").Value = 1.664315320441972

$ws.Range("B21").Value = 1.449349829459095
$ws.Range("C21").Value = 0.3496552731728713
$ws.Range("D21").Value = 0.02766443967166055
$ws.Range("F21").Value = 0.5375964350332794
$ws.Range("G21").Value = 0.3734959302236547
$ws.Range("H21").Value = 0.4969718514777668
$ws.Range("L21").Value = 0.3604284556369493
$ws.Range("N21").Value = 1.009631002091893
$ws.Range("O21").Value = 1.680530038446392

$ws.Range("B22").Value = 1.540351224298604
$ws.Range("C22").Value = 0.3535626172224937
$ws.Range("D22").Value = 0.02937285826146052
$ws.Range("F22").Value = 0.544710234560597
$ws.Range("G22").Value = 0.3789734731074645
$ws.Range("H22").Value = 0.4969046680951408
$ws.Range("L22").Value = 0.3721397415453396
$ws.Range("N22").Value = 1.007503445606659
$ws.Range("O22").Value = 1.692191215215416

$ws.Range("B23").Value = 1.491783856620998
$ws.Range("C23").Value = 0.3514750670068878
$ws.Range("D23").Value = 0.02846170397587855
$ws.Range("F23").Value = 0.5408831604719637
$ws.Range("G23").Value = 0.3760251000818329
$ws.Range("H23").Value = 0.4969188242672686
$ws.Range("L23").Value = 0.3658753415235623
$ws.Range("N23").Value = 1.008605368376891
$ws.Range("O23").Value = 1.685869964655438

$ws.Range("B24").Value = 1.307813408587947
$ws.Range("C24").Value = 0.3436165474415986
$ws.Range("D24").Value = 0.02499652022289212
$ws.Range("F24").Value = 0.5270540776907993
$ws.Range("G24").Value = 0.3654047960950351
$ws.Range("H24").Value = 0.4974491550348574
$ws.Range("L24").Value = 0.3424548864797003
$ws.Range("N24").Value = 1.013517307888876
$ws.Range("O24").Value = 1.664068235296696

$ws.Range("B25").Value = 1.109372983843514
$ws.Range("C25").Value = 0.3352483397913772
$ws.Range("D25").Value = 0.02122932807143485
$ws.Range("F25").Value = 0.5135558695824187
$ws.Range("G25").Value = 0.3551098459461315
$ws.Range("H25").Value = 0.4990394709381007
$ws.Range("L25").Value = 0.3178457721581225
$ws.Range("N25").Value = 1.020397188559372
$ws.Range("O25").Value = 1.645107730950031
